$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'328.13"
$ws.Range("E2").Value = "'0.11%"
$ws.Range("D3").Value = "'44.01"
$ws.Range("E3").Value = "'1.14%"
$ws.Range("E4").Value = "'0.55%"
$ws.Range("D5").Value = "'0.08032"
$ws.Range("E5").Value = "'-0.61%"
$ws.Range("D6").Value = "'1.972"
$ws.Range("E6").Value = "'3.99%"
$ws.Range("D7").Value = "'4.334"
$ws.Range("E7").Value = "'1.12%"
$ws.Range("D8").Value = "'0.9475"
$ws.Range("E8").Value = "'1.27%"
$ws.Range("D9").Value = "'2.552"
$ws.Range("E9").Value = "'-8.49%"
$ws.Range("D10").Value = "'0.1170"
$ws.Range("E10").Value = "'-0.84%"
$ws.Range("D11").Value = "'0.1852"
$ws.Range("E11").Value = "'-2.18%"
$ws.Range("E12").Value = "'38.34%"
$ws.Range("D13").Value = "'0.09809"
$ws.Range("E13").Value = "'2.61%"
$ws.Range("D14").Value = "'0.04720"
$ws.Range("E14").Value = "'14.30%"
$ws.Range("D15").Value = "'0.1064"
$ws.Range("E15").Value = "'-0.10%"
$ws.Range("D16").Value = "'0.001283"
$ws.Range("E16").Value = "'0.73%"
$ws.Range("D17").Value = "'0.04213"
$ws.Range("E17").Value = "'-2.78%"
$ws.Range("D18").Value = "'0.005946"
$ws.Range("E18").Value = "'0.59%"
$ws.Range("D19").Value = "'3.372"
$ws.Range("E19").Value = "'-5.53%"
$ws.Range("D20").Value = "'0.3475"
$ws.Range("E20").Value = "'-0.34%"
$ws.Range("D21").Value = "'0.1409"
$ws.Range("E21").Value = "'3.23%"
$ws.Range("D22").Value = "'0.2510"
$ws.Range("E22").Value = "'-3.07%"
$ws.Range("D23").Value = "'0.001253"
$ws.Range("E23").Value = "'1.15%"
$ws.Range("E24").Value = "'-1.82%"
$ws.Range("E25").Value = "'-2.86%"
$ws.Range("D26").Value = "'0.0003751"
$ws.Range("E26").Value = "'-6.10%"
$ws.Range("E38").Value = "'-1.95%"
$ws.Range("D39").Value = "'0.05508"
$ws.Range("E39").Value = "'0.90%"
$ws.Range("D40").Value = "'0.007549"
$ws.Range("E40").Value = "'-1.74%"
$ws.Range("E41").Value = "'0.65%"
$ws.Range("D42").Value = "'0.007649"
$ws.Range("E42").Value = "'-33.19%"
$ws.Range("D43").Value = "'0.002020"
$ws.Range("E43").Value = "'-3.99%"
$ws.Range("D44").Value = "'0.008365"
$ws.Range("E44").Value = "'-9.43%"
$ws.Range("D45").Value = "'0.00007100"
$ws.Range("E45").Value = "'1.31%"
$ws.Range("D46").Value = "'0.00000000751"
$ws.Range("E46").Value = "'0.35%"
$ws.Range("D48").Value = "'0.004842"
$ws.Range("E48").Value = "'36.17%"
$ws.Range("D49").Value = "'0.00002104"
$ws.Range("E49").Value = "'0.35%"
$ws.Range("D50").Value = "'0.0002004"
$ws.Range("E50").Value = "'0.35%"
